# Update the "From" value of rule R30 (cell C10 on sheet "Rules") from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
